$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "abc" (unused placeholder string) becomes the real "withdraw1" label for
# the D2 data cell — this also drops "abc" from the shared-string table
# (no other cell referenced it) and appends "withdraw1" as a brand-new entry.
$ws.Range("D2").Value = "withdraw1"

# Unhappy-case sample data: different account, later toDate window, and a
# much smaller observed transaction count than the configured minimum.
$ws.Range("B2").Value = 85246
$ws.Range("G2").Value = 10102020
$ws.Range("I2").Value = 20

# Column D needs to be a little wider to comfortably fit "withdraw1"/"description".
$ws.Columns("D").ColumnWidth = 10.166666666666666

# Leave the selection on G2 (toDate), matching where the edit was made.
$ws.Range("G2").Select()
